# New crime data collected - update the weekly CompStat figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + week-covering dates) ---
$ws.Range("C6").Value = "Volume 30   Number  26"
$ws.Range("C8").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# --- Plain numeric value updates ---
$ws.Range("N14").Value = -87.5
$ws.Range("N15").Value = -63.157894736842
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -47.058823529411
$ws.Range("I16").Value = 59
$ws.Range("J16").Value = 85
$ws.Range("K16").Value = -30.588235294117
$ws.Range("L16").Value = -21.333333333333
$ws.Range("M16").Value = -43.809523809523
$ws.Range("N16").Value = -86.498855835240
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 94
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 20.512820512820
$ws.Range("M17").Value = 77.358490566037
$ws.Range("N17").Value = -60.669456066945
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 157.142857142857
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 97
$ws.Range("K18").Value = -9.278350515463
$ws.Range("L18").Value = 57.142857142857
$ws.Range("M18").Value = 46.666666666666
$ws.Range("N18").Value = -83.582089552238
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -12.5
$ws.Range("I19").Value = 253
$ws.Range("J19").Value = 242
$ws.Range("K19").Value = 4.545454545454
$ws.Range("L19").Value = 29.081632653061
$ws.Range("M19").Value = -1.937984496124
$ws.Range("N19").Value = -50.294695481335
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 57
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = 54.054054054054
$ws.Range("L20").Value = 54.054054054054
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -88.690476190476
$ws.Range("C21").Value = 28
$ws.Range("E21").Value = 27.272727272727
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = 2.083333333333
$ws.Range("I21").Value = 560
$ws.Range("J21").Value = 561
$ws.Range("K21").Value = -0.178253119429
$ws.Range("L21").Value = 25
$ws.Range("M21").Value = 12
$ws.Range("N21").Value = -75.221238938053
$ws.Range("M22").Value = -50
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 54
$ws.Range("J23").Value = 68
$ws.Range("K23").Value = -20.588235294117
$ws.Range("L23").Value = -12.903225806451
$ws.Range("M23").Value = 28.571428571428
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 143
$ws.Range("H24").Value = -16.783216783216
$ws.Range("I24").Value = 889
$ws.Range("J24").Value = 899
$ws.Range("K24").Value = -1.112347052280
$ws.Range("L24").Value = 73.972602739726
$ws.Range("M24").Value = 82.921810699588
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 75
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -29.411764705882
$ws.Range("I25").Value = 152
$ws.Range("J25").Value = 155
$ws.Range("K25").Value = -1.935483870967
$ws.Range("L25").Value = 13.432835820895
$ws.Range("M25").Value = -3.184713375796
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = -6.666666666666
$ws.Range("L27").Value = 12
$ws.Range("N28").Value = -89.655172413793
$ws.Range("N29").Value = -89.655172413793
$ws.Range("G30").Value = 3
$ws.Range("J30").Value = 13
$ws.Range("K30").Value = -30.769230769230

# --- Cells that flip from a "no data" placeholder (shared text "0" / "***.*")
#     to a real number, while keeping the same number-formatted style as
#     their neighbours (copy the number format from a same-family cell). ---
$ws.Range("K14").Copy()
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("L14").Value = 100

$ws.Range("K30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1

$ws.Range("N30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100

# --- Cells that flip from a real number to the "no data" placeholder text
#     ("0" / "***.*", both already shared strings), while keeping the
#     right-aligned General-format text style used elsewhere on the sheet.
#     Force literal text via a temporary Text number format so the numeric
#     -looking string isn't re-parsed as a number, then restore the normal
#     style by pasting formats from a donor cell that already has it. ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"

$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
